$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean slate for the whole table region so stale styles
# (the old "B"/"C" highlight + left-align styles) don't leak into the
# rebuilt table.
$ws.Range("A1:H12").ClearContents()
$ws.Range("A1:H12").ClearFormats()

# Header row
$ws.Range("A1").Value = "nome "
$ws.Range("B1").Value = "empresa "
$ws.Range("C1").Value = "telefone"

# Two plain (unstyled) rows
$ws.Range("A2").Value = "dani"
$ws.Range("B2").Value = "grupar"
$ws.Range("C2").Value = "55 44 9101-8419"

$ws.Range("A3").Value = "gabriel "
$ws.Range("B3").Value = "grupar"
$ws.Range("C3").Value = "55 44 9871-6404"

# Remaining rows (center-aligned name/company columns)
$data = @(
    @("Gisele",   "BAS AGROFLORESTAL", "55 38 9981-8469"),
    @("Eric",      "MILOG",             "55 44 9928-1181"),
    @("Eric",      "MILOG",             "55 44 9865-9762"),
    @("Nelis",     "AVILA",             "55 55 8142-0373"),
    @("Adrina",    "BETTENCOURT",       "55 45 9965-7343"),
    @("Thiago",    "THIAGO",            "55 14 9678-0046"),
    @("Stark",     "STARK METAIS",      "55 44 9901-3774"),
    @("Anderson",  "NARDOCI MOCCHI",    " 55 44 9986-0190"),
    @("Vinicius",  "NEDI",              "55 53 8425-9087")
)

$r = 4
foreach ($row in $data) {
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $r++
}

$ws.Range("A4:B12").HorizontalAlignment = -4108

# Stray formatted-but-empty cell left over at H11 (matches the highlighted
# "grupar" style used before: white fill, black font, left aligned)
$h11 = $ws.Range("H11")
$h11.HorizontalAlignment = -4131
$h11.Font.Color = 0
$h11.Interior.Pattern = 1
$h11.Interior.Color = 16777215

$ws.Range("H9").Select()
